# Emails.xlsx - "Criado disparador de e-mail, e corrigido arquivo de e-mail"
#
# The e-mail addresses used for the dispatcher test accounts moved from the
# "pythonimpressionador+<nome>@gmail.com" alias scheme to
# "romuloptmota+<nome>@gmail.com", and the old per-cell mailto hyperlinks
# (with their relationship parts) are dropped in favour of plain text in
# column C. The last-used cell/view position also moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $manager = [string]$ws.Cells.Item($r, 2).Value2
    if ([string]::IsNullOrEmpty($manager)) { continue }

    $slug = $manager.ToLower().Replace(" ", "_")
    $ws.Cells.Item($r, 3).Value2 = "romuloptmota+$slug@gmail.com"
}

# The refreshed workbook no longer carries the mailto hyperlinks on column C
# (plain text keeps the existing "Hiperlink" cell style).
[void]$ws.Hyperlinks.Delete()

# Last selection left on the sheet when it was saved.
[void]$ws.Range("G20").Select()

$wb.Save()
